$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns that would otherwise be auto-detected as numbers to stay as Text,
# matching the source data (which stores these as literal strings).
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D15",
    "D16",
    "D20",
    "D21",
    "D22",
    "D23",
    "D26",
    "D27",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D40",
    "D41",
    "D42",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.134.94"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "2.057.22"
$ws.Range("E3").Value = "  -3.26%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "248.85"
$ws.Range("E5").Value = "  -2.91%  "

$ws.Range("D6").Value = "0.655"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "54.82"
$ws.Range("E8").Value = "  +16.23%  "

$ws.Range("D9").Value = "61.99"
$ws.Range("E9").Value = "  +3.65%  "

$ws.Range("D10").Value = "0.378"
$ws.Range("E10").Value = "  +0.97%  "

$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  +5.23%  "

$ws.Range("E12").Value = "  +5.39%  "

$ws.Range("D13").Value = "15.09"
$ws.Range("E13").Value = "  +5.07%  "

$ws.Range("D14").Value = "2.355.33"
$ws.Range("E14").Value = "  -3.34%  "

$ws.Range("D15").Value = "0.819"
$ws.Range("E15").Value = "  -2.58%  "

$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "2.054.96"
$ws.Range("E17").Value = "  -3.34%  "

$ws.Range("D18").Value = "37.073.48"
$ws.Range("E18").Value = "  +0.92%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0905"
$ws.Range("E19").Value = "  +7.87%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "72.37"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  +7.90%  "

$ws.Range("D22").Value = "5.31"
$ws.Range("E22").Value = "  +1.81%  "

$ws.Range("D23").Value = "237.05"
$ws.Range("E23").Value = "  -2.00%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("D26").Value = "169.97"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -2.14%  "

$ws.Range("E28").Value = "  -7.98%  "

$ws.Range("E29").Value = "  -3.15%  "

$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").Value = "4.56"
$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("E32").Value = "  +11.29%  "

$ws.Range("D33").Value = "0.0624"
$ws.Range("E33").Value = "  +3.85%  "

$ws.Range("D34").Value = "4.32"
$ws.Range("E34").Value = "  +3.61%  "

$ws.Range("D35").Value = "0.0880"
$ws.Range("E35").Value = "  -8.72%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  -4.89%  "

$ws.Range("E38").Value = "  -7.93%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  +20.77%  "

$ws.Range("D41").Value = "18.29"
$ws.Range("E41").Value = "  +12.67%  "

$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").Value = "4.42"
$ws.Range("E44").Value = "  +48.49%  "

$ws.Range("D45").Value = "96.06"
$ws.Range("E45").Value = "  -3.24%  "

$ws.Range("B46").Value = "Gas"
$ws.Range("C46").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D46").Value = "14.70"
$ws.Range("E46").Value = "  -49.88%  "

$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +5.17%  "

$ws.Range("D49").Value = "1.297.52"
$ws.Range("E49").Value = "  -4.61%  "

$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").Value = "6.79"
$ws.Range("E51").Value = "  -6.54%  "
